# edit.ps1 - Applies the CSFG7.pptx "lecture 7" update:
#   Slide 11 ("CORRELATION HYPOTHESIS TESTS"): the Spearman p-value in the
#   line "... Spearman R is 0.19, p=0.067 - not significant!" is corrected
#   from 0.067 to 0.087.
#
# (The accompanying diff also shows slide 7's "Some tests are two-tailed
# tests" text box being re-saved with refreshed dirty="0" proofing flags
# but no textual/formatting change; that is a cosmetic no-op left by the
# original PowerPoint editor and is not reproduced here.)

function Get-ShapeById {
    param($Slide, [int]$Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# --- Slide 11: "... Spearman R is 0.19, p=0.067 - not significant!" ---
#     Update the p-value from 0.067 to 0.087 (keeping the trailing space
#     that precedes the en-dash, matching how the run ends up split in
#     the saved file).
$slide11 = $p.Slides.Item(11)
$shape11 = Get-ShapeById $slide11 11

$tr11 = $shape11.TextFrame.TextRange
$bodyText = $tr11.Text
$oldValue = "0.067 "
$newValue = "0.087 "
$pos = $bodyText.IndexOf($oldValue)
if ($pos -ge 0) {
    $run = $tr11.Characters($pos + 1, $oldValue.Length)
    $run.Text = $newValue
}
